# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement detail table (rows 16-31, cols C:F) is rebuilt:
# the old two employees (JAIDER DARIO ZAPATEIRO SALGADO / DAVID JOSE RUIZ
# PUERTA) were interleaved period-by-period; the new data groups all
# periods for DAVID JOSE RUIZ PUERTA first (1907 down to 1812), followed
# by all periods for JAIDER DARIO ZAPATEIRO SALGADO (1907 down to 1812).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; Doc = "1047456907"; Name = "DAVID JOSE RUIZ PUERTA";           Period = "1907"; Value = 15625 },
    @{ Row = 17; Doc = "1047456907"; Name = "DAVID JOSE RUIZ PUERTA";           Period = "1906"; Value = 31249 },
    @{ Row = 18; Doc = "1047456907"; Name = "DAVID JOSE RUIZ PUERTA";           Period = "1905"; Value = 31249 },
    @{ Row = 19; Doc = "1047456907"; Name = "DAVID JOSE RUIZ PUERTA";           Period = "1904"; Value = 31249 },
    @{ Row = 20; Doc = "1047456907"; Name = "DAVID JOSE RUIZ PUERTA";           Period = "1903"; Value = 31249 },
    @{ Row = 21; Doc = "1047456907"; Name = "DAVID JOSE RUIZ PUERTA";           Period = "1902"; Value = 31249 },
    @{ Row = 22; Doc = "1047456907"; Name = "DAVID JOSE RUIZ PUERTA";           Period = "1901"; Value = 31249 },
    @{ Row = 23; Doc = "1047456907"; Name = "DAVID JOSE RUIZ PUERTA";           Period = "1812"; Value = 31249 },
    @{ Row = 24; Doc = "1143387469"; Name = "JAIDER DARIO ZAPATEIRO SALGADO";   Period = "1907"; Value = 15625 },
    @{ Row = 25; Doc = "1143387469"; Name = "JAIDER DARIO ZAPATEIRO SALGADO";   Period = "1906"; Value = 31249 },
    @{ Row = 26; Doc = "1143387469"; Name = "JAIDER DARIO ZAPATEIRO SALGADO";   Period = "1905"; Value = 31249 },
    @{ Row = 27; Doc = "1143387469"; Name = "JAIDER DARIO ZAPATEIRO SALGADO";   Period = "1904"; Value = 31249 },
    @{ Row = 28; Doc = "1143387469"; Name = "JAIDER DARIO ZAPATEIRO SALGADO";   Period = "1903"; Value = 31249 },
    @{ Row = 29; Doc = "1143387469"; Name = "JAIDER DARIO ZAPATEIRO SALGADO";   Period = "1902"; Value = 31249 },
    @{ Row = 30; Doc = "1143387469"; Name = "JAIDER DARIO ZAPATEIRO SALGADO";   Period = "1901"; Value = 31249 },
    @{ Row = 31; Doc = "1143387469"; Name = "JAIDER DARIO ZAPATEIRO SALGADO";   Period = "1812"; Value = 31249 }
)

foreach ($r in $rows) {
    $ws.Range("C" + $r.Row).Value = $r.Doc
    $ws.Range("D" + $r.Row).Value = $r.Name
    $ws.Range("E" + $r.Row).Value = $r.Period
    $ws.Range("F" + $r.Row).Value = $r.Value
}
